# Update "想去人数" (number of people interested) figures for a handful of
# events on both the "展览" sheet and the aggregated "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    6  = 261
    7  = 6304
    10 = 105
    11 = 67
    15 = 475
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
